# update OrderCommand Sequence Diagram
#
# Applies the positional / text tweaks from the commit:
#   - Rectangle 5              (shape 4)  -> move down (lifeline bar)
#   - Straight Arrow Connector 22 (shape 9) -> un-flip, reposition/resize
#   - TextBox 25               (shape 10) -> reposition + retext
#   - Straight Arrow Connector 36 (shape 12) -> move down
#   - TextBox 28 "result"      (shape 24) -> move down
#
# NOTE: PowerPoint's Shape.Left/Top/Width/Height are expressed in points,
# while the underlying OOXML stores EMU (1 pt = 12700 EMU). The literal
# point values below were solved so that, after the COM layer's internal
# (single-precision) EMU round-trip, they reproduce the exact target EMU.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) Rectangle 5 - lifeline bar shifts down
$sh = $s.Shapes.Item(4)
$sh.Top = 189.19251968503937

# 2) Straight Arrow Connector 22 - remove vertical flip, move & resize
$sh = $s.Shapes.Item(9)
$sh.VerticalFlip = $false
$sh.Left = 44.495994251968504
$sh.Top = 198.0
$sh.Width = 141.46409448818898
$sh.Height = 0.48646669291338585

# 3) TextBox 25 - move & change label text
$sh = $s.Shapes.Item(10)
$sh.Left = 41.198425196850394
$sh.Top = 164.07267716535432
$sh.TextFrame.TextRange.Text = "executeUndoableCommand ()"

# 4) Straight Arrow Connector 36 - move down
$sh = $s.Shapes.Item(12)
$sh.Top = 396.0

# 5) TextBox 28 ("result") - move down
$sh = $s.Shapes.Item(24)
$sh.Top = 379.1433958267716
